$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.77%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.115"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.21%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07621"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.63%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.605"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.41%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9029"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.37%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1111"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8.27%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1785"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.08%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09143"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.14%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04203"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.41%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1052"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.39%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001249"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.59%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005670"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.55%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.347"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.25%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.246"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.35%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.52%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.582"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.09%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1365"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.89%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.47%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04066"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.67%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001230"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.33%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004123"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001300"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.15%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003746"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.70%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05188"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.58%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007768"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.26%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1302"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.95%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007046"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "10.89%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001951"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.47%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008782"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.26%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3334"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.32%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006937"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.77%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03160"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "605.19%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.04%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.04%"
